# Auto-generated edit script: updates cryptocurrency Price (col D)
# and Volume(1h) (col E) values on Sheet1, matching the commit diff.
# A leading apostrophe forces Excel to treat numeric-looking values
# (e.g. "591.96") as text, matching the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'64.382.66"
$ws.Cells.Item(3, 4).Value = "'3.509.36"
$ws.Cells.Item(3, 5).Value = "'  +0.48%  "
$ws.Cells.Item(4, 5).Value = "'  -0.01%  "
$ws.Cells.Item(5, 4).Value = "'591.96"
$ws.Cells.Item(5, 5).Value = "'  +0.88%  "
$ws.Cells.Item(6, 4).Value = "'134.68"
$ws.Cells.Item(6, 5).Value = "'  -0.13%  "
$ws.Cells.Item(8, 5).Value = "'  +0.21%  "
$ws.Cells.Item(9, 5).Value = "'  +5.77%  "
$ws.Cells.Item(10, 5).Value = "'  +0.67%  "
$ws.Cells.Item(11, 4).Value = "'0.388"
$ws.Cells.Item(11, 5).Value = "'  +3.33%  "
$ws.Cells.Item(12, 4).Value = "'4.107.76"
$ws.Cells.Item(12, 5).Value = "'  +0.37%  "
$ws.Cells.Item(13, 5).Value = "'  +0.98%  "
$ws.Cells.Item(14, 5).Value = "'  +0.56%  "
$ws.Cells.Item(15, 4).Value = "'3.511.22"
$ws.Cells.Item(15, 5).Value = "'  +0.27%  "
$ws.Cells.Item(16, 4).Value = "'25.74"
$ws.Cells.Item(16, 5).Value = "'  +1.67%  "
$ws.Cells.Item(17, 4).Value = "'64.378.17"
$ws.Cells.Item(17, 5).Value = "'  +0.13%  "
$ws.Cells.Item(18, 4).Value = "'9.97"
$ws.Cells.Item(18, 5).Value = "'  -0.36%  "
$ws.Cells.Item(19, 4).Value = "'13.64"
$ws.Cells.Item(19, 5).Value = "'  -0.89%  "
$ws.Cells.Item(20, 4).Value = "'5.75"
$ws.Cells.Item(20, 5).Value = "'  +2.01%  "
$ws.Cells.Item(21, 4).Value = "'390.40"
$ws.Cells.Item(21, 5).Value = "'  +1.17%  "
$ws.Cells.Item(22, 5).Value = "'  +2.19%  "
$ws.Cells.Item(23, 4).Value = "'3.649.45"
$ws.Cells.Item(23, 5).Value = "'  +0.43%  "
$ws.Cells.Item(24, 4).Value = "'74.59"
$ws.Cells.Item(24, 5).Value = "'  +0.69%  "
$ws.Cells.Item(25, 5).Value = "'  +0.05%  "
$ws.Cells.Item(26, 4).Value = "'5.75"
$ws.Cells.Item(26, 5).Value = "'  +0.54%  "
$ws.Cells.Item(27, 5).Value = "'  +3.27%  "
$ws.Cells.Item(28, 5).Value = "'  +0.02%  "
$ws.Cells.Item(29, 5).Value = "'  -0.50%  "
$ws.Cells.Item(30, 5).Value = "'  +1.51%  "
$ws.Cells.Item(31, 5).Value = "'  +0.47%  "
$ws.Cells.Item(32, 4).Value = "'1.47"
$ws.Cells.Item(32, 5).Value = "'  -4.53%  "
$ws.Cells.Item(33, 5).Value = "'  +6.99%  "
$ws.Cells.Item(34, 4).Value = "'3.535.41"
$ws.Cells.Item(34, 5).Value = "'  +0.57%  "
$ws.Cells.Item(35, 5).Value = "'  +0.00%  "
$ws.Cells.Item(36, 4).Value = "'23.36"
$ws.Cells.Item(36, 5).Value = "'  -0.37%  "
$ws.Cells.Item(37, 4).Value = "'5.37"
$ws.Cells.Item(37, 5).Value = "'  +1.39%  "
$ws.Cells.Item(38, 4).Value = "'6.95"
$ws.Cells.Item(38, 5).Value = "'  +1.68%  "
$ws.Cells.Item(39, 5).Value = "'  +1.23%  "
$ws.Cells.Item(40, 4).Value = "'166.86"
$ws.Cells.Item(40, 5).Value = "'  +2.69%  "
$ws.Cells.Item(41, 5).Value = "'  +0.91%  "
$ws.Cells.Item(42, 5).Value = "'  +1.03%  "
$ws.Cells.Item(43, 5).Value = "'  +0.00%  "
$ws.Cells.Item(44, 5).Value = "'  +0.94%  "
$ws.Cells.Item(45, 4).Value = "'24.86"
$ws.Cells.Item(45, 5).Value = "'  -3.92%  "
$ws.Cells.Item(46, 5).Value = "'  +0.28%  "
$ws.Cells.Item(47, 5).Value = "'  -2.93%  "
$ws.Cells.Item(48, 4).Value = "'6.80"
$ws.Cells.Item(48, 5).Value = "'  +0.70%  "
$ws.Cells.Item(49, 5).Value = "'  +0.73%  "
$ws.Cells.Item(50, 4).Value = "'2.373.05"
$ws.Cells.Item(50, 5).Value = "'  -4.17%  "
$ws.Cells.Item(51, 5).Value = "'  +0.43%  "
